# Appends the new game item rows (43-49) to the storage sheet, replacing the
# previously duplicated/repeated hard-coded item blocks with a single
# data-driven loop that writes each record into the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$items = @(
    @("BonusPower", 0.11, 0,    49,  100, 0,   "lose"),
    @("SkipBoss",   0.04, 0,    25,  20,  0,   "lose"),
    @("SkipBoss",   0.03, 0,    22,  200, 0,   "lose"),
    @("BonusPower", 2,    -40,  156, 100, 200, "win"),
    @("SkipBoss",   0.04, 0,    21,  20,  0,   "lose"),
    @("SkipBoss",   2,    500,  132, 20,  40,  "win"),
    @("BonusPower", 2,    810,  131, 100, 200, "win")
)

$startRow = 43
for ($i = 0; $i -lt $items.Count; $i++) {
    $row = $startRow + $i
    $record = $items[$i]

    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
}
